# Add a second slide to the deck (commit: "add a second slide to test
# scrolling"). The new slide uses the "Title and Content" layout
# (CustomLayout #2 on the slide master, i.e. the legacy ppLayoutText = 2),
# matching ppt/slideLayouts/slideLayout2.xml, and is appended after the
# existing slide.

$p = $ppt.ActivePresentation

# Append a new slide (index 2 = after the current last slide) using the
# "Title and Content" layout.
$slide = $p.Slides.Add(2, 2)

# Title placeholder -> "Objectives" (with trailing tab, as authored).
$title = $slide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Objectives`t"

# Body / content placeholder -> six bullet lines, one per paragraph.
$bullets = @(
    "Know what a Linux terminal is.",
    "Know how to access a Linux terminal.",
    "Know the most important Linux terminal commands.",
    "Install Python",
    "Write a “Hello World” Python script",
    "Run the “Hello World” Python script"
)
$body = $slide.Shapes.Item(2).TextFrame.TextRange
$body.Text = [string]::Join("`r", $bullets)
